$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 165.64706
$ws.Range("I11").Value = 165.64706
$ws.Range("K11").Value = 165.64706
$ws.Range("M11").Value = -25.64706000000001
$ws.Range("H28").Value = 2820.6
$ws.Range("I28").Value = 1049
$ws.Range("J28").Value = 4001.6667
$ws.Range("K28").Value = 1049
$ws.Range("L28").Value = 4001.6667
$ws.Range("M28").Value = -564
$ws.Range("N28").Value = -4971.6667
$ws.Range("H29").Value = 2887.5
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15562
$ws.Range("H37").Value = 914.75
$ws.Range("I37").Value = 919.6667
$ws.Range("K37").Value = 2759.0001
$ws.Range("M37").Value = -2633.0001
$ws.Range("H38").Value = 18443.834
$ws.Range("I38").Value = 25166
$ws.Range("J38").Value = 4999.5
$ws.Range("K38").Value = 75498
$ws.Range("L38").Value = 14998.5
$ws.Range("M38").Value = -75126
$ws.Range("N38").Value = -15742.5
$ws.Range("H107").Value = 359.2143
$ws.Range("I107").Value = 366.23077
$ws.Range("K107").Value = 366.23077
$ws.Range("M107").Value = 1553.76923
$ws.Range("H116").Value = 5942.6665
$ws.Range("I116").Value = 5580.6665
$ws.Range("K116").Value = 5580.6665
$ws.Range("M116").Value = -2138.6665
$ws.Range("H137").Value = 1391337.5
$ws.Range("I137").Value = 1924595.8
$ws.Range("J137").Value = 4866
$ws.Range("K137").Value = 5773787.4
$ws.Range("L137").Value = 14598
$ws.Range("M137").Value = -5771237.4
$ws.Range("N137").Value = -19698
$ws.Range("H141").Value = 1945.3334
$ws.Range("J141").Value = 932.6667
$ws.Range("L141").Value = 2798.0001
$ws.Range("N141").Value = -13158.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2052.5293
$ws.Range("I32").Value = 2100.0303
$ws.Range("K32").Value = 2100.0303
$ws.Range("M32").Value = -1813.0303
$ws.Range("H63").Value = 4342.7144
$ws.Range("I63").Value = 2224.75
$ws.Range("J63").Value = 7166.6665
$ws.Range("K63").Value = 2224.75
$ws.Range("L63").Value = 7166.6665
$ws.Range("M63").Value = -1538.75
$ws.Range("N63").Value = -8538.666499999999
$ws.Range("H66").Value = 4342.7144
$ws.Range("I66").Value = 2224.75
$ws.Range("J66").Value = 7166.6665
$ws.Range("K66").Value = 11123.75
$ws.Range("L66").Value = 35833.3325
$ws.Range("M66").Value = -7691.75
$ws.Range("N66").Value = -42697.3325
$ws.Range("H74").Value = 3630.7144
$ws.Range("I74").Value = 3232.5356
$ws.Range("J74").Value = 5223.4287
$ws.Range("K74").Value = 3232.5356
$ws.Range("L74").Value = 5223.4287
$ws.Range("M74").Value = -2358.5356
$ws.Range("N74").Value = -6971.4287
$ws.Range("H77").Value = 3630.7144
$ws.Range("I77").Value = 3232.5356
$ws.Range("J77").Value = 5223.4287
$ws.Range("K77").Value = 16162.678
$ws.Range("L77").Value = 26117.1435
$ws.Range("M77").Value = -11794.678
$ws.Range("N77").Value = -34853.14350000001
$ws.Range("H88").Value = 5815.8335
$ws.Range("I88").Value = 4000
$ws.Range("J88").Value = 6179
$ws.Range("K88").Value = 4000
$ws.Range("L88").Value = 6179
$ws.Range("M88").Value = -3594
$ws.Range("N88").Value = -6991
$ws.Range("H91").Value = 5815.8335
$ws.Range("I91").Value = 4000
$ws.Range("J91").Value = 6179
$ws.Range("K91").Value = 4000
$ws.Range("L91").Value = 6179
$ws.Range("M91").Value = -2596
$ws.Range("N91").Value = -8987
$ws.Range("H110").Value = 3412.1177
$ws.Range("J110").Value = 1722
$ws.Range("L110").Value = 1722
$ws.Range("N110").Value = -5812
$ws.Range("H132").Value = 2835
$ws.Range("I132").Value = 1687.3
$ws.Range("K132").Value = 5061.9
$ws.Range("M132").Value = -2531.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 6999.8
$ws.Range("J36").Value = 3333
$ws.Range("L36").Value = 3333
$ws.Range("N36").Value = -4401
$ws.Range("H99").Value = 1487.7
$ws.Range("J99").Value = 1798.1
$ws.Range("L99").Value = 1798.1
$ws.Range("N99").Value = -4794.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2221.6
$ws.Range("I16").Value = 2221
$ws.Range("J16").Value = 2222
$ws.Range("K16").Value = 2221
$ws.Range("L16").Value = 2222
$ws.Range("M16").Value = -1934
$ws.Range("N16").Value = -2796
$ws.Range("H31").Value = 4189.3213
$ws.Range("I31").Value = 2932.647
$ws.Range("K31").Value = 2932.647
$ws.Range("M31").Value = -2637.647
$ws.Range("H34").Value = 4189.3213
$ws.Range("I34").Value = 2932.647
$ws.Range("K34").Value = 2932.647
$ws.Range("M34").Value = -2730.647
$ws.Range("H62").Value = 7695208
$ws.Range("I62").Value = 11113816
$ws.Range("K62").Value = 11113816
$ws.Range("M62").Value = -11113192
$ws.Range("H65").Value = 7695208
$ws.Range("I65").Value = 11113816
$ws.Range("K65").Value = 55569080
$ws.Range("M65").Value = -55565960
$ws.Range("H99").Value = 5141.4287
$ws.Range("I99").Value = 1506
$ws.Range("J99").Value = 6595.6
$ws.Range("K99").Value = 1506
$ws.Range("L99").Value = 6595.6
$ws.Range("M99").Value = -8
$ws.Range("N99").Value = -9591.6
$ws.Range("H105").Value = 1787.3636
$ws.Range("I105").Value = 1381.5
$ws.Range("J105").Value = 2869.6667
$ws.Range("K105").Value = 1381.5
$ws.Range("L105").Value = 2869.6667
$ws.Range("M105").Value = 365.5
$ws.Range("N105").Value = -6363.6667
$ws.Range("H113").Value = 2221.6
$ws.Range("I113").Value = 2221
$ws.Range("J113").Value = 2222
$ws.Range("K113").Value = 2221
$ws.Range("L113").Value = 2222
$ws.Range("M113").Value = -51
$ws.Range("N113").Value = -6562
$ws.Range("H126").Value = 5141.4287
$ws.Range("I126").Value = 1506
$ws.Range("J126").Value = 6595.6
$ws.Range("K126").Value = 4518
$ws.Range("L126").Value = 19786.8
$ws.Range("M126").Value = -2048
$ws.Range("N126").Value = -24726.8
$ws.Range("H134").Value = 2629.1177
$ws.Range("I134").Value = 2592.08
$ws.Range("K134").Value = 7776.24
$ws.Range("M134").Value = -5241.24

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 9316.556
$ws.Range("I23").Value = 417.5
$ws.Range("J23").Value = 16435.8
$ws.Range("K23").Value = 1252.5
$ws.Range("L23").Value = 49307.39999999999
$ws.Range("M23").Value = -1017.5
$ws.Range("N23").Value = -49777.39999999999
$ws.Range("H114").Value = 3759.6
$ws.Range("I114").Value = 1932.6666
$ws.Range("J114").Value = 6500
$ws.Range("K114").Value = 5797.9998
$ws.Range("L114").Value = 19500
$ws.Range("M114").Value = -2543.9998
$ws.Range("N114").Value = -26008
$ws.Range("H117").Value = 1252
$ws.Range("J117").Value = 1144.6666
$ws.Range("L117").Value = 3433.9998
$ws.Range("N117").Value = -10317.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1980.6
$ws.Range("I102").Value = 829.6667
$ws.Range("K102").Value = 829.6667
$ws.Range("M102").Value = 792.3333
$ws.Range("H103").Value = 89999
$ws.Range("J103").Value = 89999
$ws.Range("L103").Value = 89999
$ws.Range("N103").Value = -92343
$ws.Range("H132").Value = 2709.7407
$ws.Range("I132").Value = 2113.9473
$ws.Range("J132").Value = 4124.75
$ws.Range("K132").Value = 6341.841899999999
$ws.Range("L132").Value = 12374.25
$ws.Range("M132").Value = -3811.841899999999
$ws.Range("N132").Value = -17434.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5827
$ws.Range("I7").Value = 4556.1665
$ws.Range("K7").Value = 4556.1665
$ws.Range("M7").Value = -4444.1665
$ws.Range("H61").Value = 3592.2
$ws.Range("I61").Value = 3548.5557
$ws.Range("K61").Value = 3548.5557
$ws.Range("M61").Value = -3346.5557
$ws.Range("H113").Value = 3592.2
$ws.Range("I113").Value = 3548.5557
$ws.Range("K113").Value = 3548.5557
$ws.Range("M113").Value = -1378.5557
$ws.Range("H126").Value = 5827
$ws.Range("I126").Value = 4556.1665
$ws.Range("K126").Value = 13668.4995
$ws.Range("M126").Value = -11198.4995
$ws.Range("H132").Value = 4553.579
$ws.Range("I132").Value = 3736.0688
$ws.Range("K132").Value = 11208.2064
$ws.Range("M132").Value = -8678.206399999999
$ws.Range("H136").Value = 3365.182
$ws.Range("I136").Value = 3029.7144
$ws.Range("J136").Value = 3952.25
$ws.Range("K136").Value = 9089.143199999999
$ws.Range("L136").Value = 11856.75
$ws.Range("M136").Value = -6539.143199999999
$ws.Range("N136").Value = -16956.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2276.3
$ws.Range("I2").Value = 2276.3
$ws.Range("K2").Value = 2276.3
$ws.Range("M2").Value = -2164.3
$ws.Range("H4").Value = 1357.5
$ws.Range("I4").Value = 2082
$ws.Range("K4").Value = 2082
$ws.Range("M4").Value = -1969
$ws.Range("H21").Value = 20007.5
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20470
$ws.Range("H35").Value = 20007.5
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20580
$ws.Range("H132").Value = 3384.279
$ws.Range("I132").Value = 2808
$ws.Range("K132").Value = 8424
$ws.Range("M132").Value = -5894
